$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the content of row 11 and row 12 for the columns that differ
# between the two rows (A, B, D, E, F, G, H, Q, R). The other columns
# (C, I, K, P, S, T, U, V, W, Y, Z, AA, AB, AD, AE, AG, AT, AW, AX, AY)
# hold identical values in both rows so they are left untouched.

$columns = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

foreach ($col in $columns) {
    $cell11 = $ws.Range($col + "11")
    $cell12 = $ws.Range($col + "12")

    $value11 = $cell11.Value2
    $value12 = $cell12.Value2

    $cell11.Value2 = $value12
    $cell12.Value2 = $value11
}
